$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column S, row 4: header value "2022", same formatting as R4 (year header)
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

# New column S, row 5: data value 42, formatted like R5 but with a "0.0" number format
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 42
$ws.Range("S5").NumberFormat = "0.0"

$excel.CutCopyMode = $false

# Update the active selection to match the new state of the workbook
$ws.Range("U4").Select()
